# Generate Report for Handback
# Updates the zh-cn / de-de handback status rows: mark them handed back
# (in sync with en-US), refresh the "Latest Handback DateTime" stamps,
# and clear the stale "version not latest" error detail. Also widens the
# Status / Error Detail columns so the new text fits.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# --- Overview sheet, data row 2 (status summary columns mirror the per-locale Status) ---
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet, data row 2 ---
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-10-18 04:16:33"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet, data row 2 ---
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-10-18 04:16:57"
$wsDeDe.Range("P2").Value = ""

# --- column width adjustments to fit the new status/error text ---
$wsOverview.Columns(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns(16).ColumnWidth = 12.833333333333334
